# Insert a new data row into the "Apio" price sheet right before the existing
# row 324, shifting all following rows down by one (dimension grows from
# A1:R445 to A1:R446), then populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a blank row at position 324 (pushes old rows 324..445 to 325..446)
$ws.Rows.Item(324).Insert()

# Populate the newly inserted row 324 with the new record
$ws.Range("A324").Value = 3
$ws.Range("B324").Value = "Femacal de La Calera"
$ws.Range("C324").Value = "Coquimbo"
$ws.Range("D324").Value = 44795
$ws.Range("E324").Value = 5
$ws.Range("F324").Value = 100112017
$ws.Range("G324").Value = "Apio"
$ws.Range("H324").Value = "Americana (o)"
$ws.Range("I324").Value = "Primera"
$ws.Range("J324").Value = 230
$ws.Range("K324").Value = 9500
$ws.Range("L324").Value = 10000
$ws.Range("M324").Value = 9761
$ws.Range("N324").Value = "$/docena de matas"
$ws.Range("O324").Value = "Pan de Azúcar"
$ws.Range("P324").Value = 1627
$ws.Range("Q324").Value = 6
$ws.Range("R324").Value = "Hortaliza"
